$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that directly follows the
#    "Play Dream Drop Diamonds Free | Slot Game Review" title paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Text -match "Meta description") {
    $p.Range.Delete()
    break
  }
}

# 2. Replace the final paragraph (the italic "Create an eye-catching
#    feature image..." image-prompt paragraph) with two paragraphs:
#      a) a bold "Play Dream Drop Diamonds Free | Slot Game Review" line
#      b) an italic meta-description line with the new copy text.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$fullRange = $lastPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dream Drop Diamonds Free | Slot Game Review</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Get ready to play Dream Drop Diamonds for free online! Read our review to learn more about its features, graphics, and jackpots.</w:t></w:r></w:p>'

$fullRange.InsertXML($xml) | Out-Null
